$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "baivab"
$ws.Range("B4").Value = "nayak"
$ws.Range("C4").Value = "baivab@gmail.com"

$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:baivab@gmail.com")

$ws.Range("A4").Select()
